# DEV: Start refactoring export category
#
# Renames the "export" category's service rows (category = "export",
# D column 7030-7035) from the old "orders"-centric naming to the new
# generic "data" naming:
#   export-data-orders-pdf   -> export-data-pdf
#   export-data-orders-csv   -> export-data-csv
#   export-data-orders-xlsx  -> export-data-xlsx
#   export-data-orders-json  -> export-data-json
#   export-data-orders-xml   -> export-data-xml
#   export-manager           -> export-data-manager
# and the matching base-path (column E) values move from
# /orderizer/export/orders/<fmt>/ to /orderizer/export/data/<fmt>/v1/
# (and /orderizer/export/manager/ -> /orderizer/export/data/manager/v1/).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 - pdf
$ws.Range("C23").Value = "export-data-pdf"
$ws.Range("E23").Value = "/orderizer/export/data/pdf/v1/"

# Row 24 - csv
$ws.Range("C24").Value = "export-data-csv"
$ws.Range("E24").Value = "/orderizer/export/data/csv/v1/"

# Row 25 - xlsx
$ws.Range("C25").Value = "export-data-xlsx"
$ws.Range("E25").Value = "/orderizer/export/data/xlsx/v1/"

# Row 26 - json
$ws.Range("C26").Value = "export-data-json"
$ws.Range("E26").Value = "/orderizer/export/data/json/v1/"

# Row 27 - xml
$ws.Range("C27").Value = "export-data-xml"
$ws.Range("E27").Value = "/orderizer/export/data/xml/v1/"

# Row 28 - manager
$ws.Range("C28").Value = "export-data-manager"
$ws.Range("E28").Value = "/orderizer/export/data/manager/v1/"

# Update the cursor/selection left where the author left it while editing.
$ws.Range("A10").Select()
